# Dietician_testdata.xlsx update — Userlogin DDT data refresh (negative/positive rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing admin login moves from column B/A swap-free update ---
# A2 keeps the "Team9.admin@gmail.com" mailto hyperlink but the cell text changes
$ws.Range("A2").Value = "Team9.admin@gmail.com"
$ws.Range("B2").Value = "test"

# --- Row 3: new negative-test login/password pair ---
$ws.Range("A3").Value = "testing@gmail.com"
$ws.Range("B3").Value = "sample"

# --- New column C: scenario / login-id labels ---
$ws.Range("C1").Value = "scenario"
$ws.Range("C2").Value = "Login1"
$ws.Range("C3").Value = "Login2"

# Style the new C3 cell with a distinct font (10pt blue Courier New)
$c3Font = $ws.Range("C3").Font
$c3Font.Name = "Courier New"
$c3Font.Size = 10
$c3Font.Color = 16744448

# Leave the selection on the newly edited cell, matching the saved workbook state
$ws.Range("C3").Select()
